$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

$ws.Range("D6").Value = "дневные по графику"
$ws.Range("D7").Value = "ночные по графику"
$ws.Range("D8").Value = "отхожено дневных"
$ws.Range("D9").Value = "отхожено ночных"

$ws.Range("E5").Value = "дни"
$ws.Range("F5").Value = "часы"
$ws.Range("G5").Value = "проценты"

$ws.Range("E6").Value = 6
$ws.Range("E8").Value = 6
$ws.Range("E9").Value = 6

$ws.Range("F12").Formula = "=E8*11.7"
$ws.Range("F13").Formula = "=E9*11"

$ws.Range("J12").Formula = "=I12*F12"
$ws.Range("J13").Formula = "=I13*F13"

$ws.Range("L13").Formula = "=J12+J13"

$ws.Range("L15").Select()
